$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Percentage-formatted cells need NumberFormat forced to Text ("@") before
# assignment, otherwise Excel auto-converts strings like "80%" into the
# numeric value 0.8 with a Percentage format (exactly like typing it in the UI).
$percentCells = @(
    "H2",
    "H5",
    "H9",
    "H11",
    "H12",
    "H14",
    "H15",
    "H16",
    "H17",
    "H19",
    "H21",
    "H23",
    "H24",
    "H25",
    "H27",
    "H28",
    "H30",
    "H31",
    "H32",
    "H34",
    "H36"
)
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-07 18:47:53"
$ws.Range("H2").Value = "80%"
$ws.Range("K2").Value = "8.6 MJ/m2"
$ws.Range("O2").Value = "-0.6 °C"
$ws.Range("E3").Value = "2026-02-07 18:47:56"
$ws.Range("K3").Value = "14.3 MJ/m2"
$ws.Range("L3").Value = "33.5 km/h - 116º 17:27 TU"
$ws.Range("E4").Value = "2026-02-07 18:47:58"
$ws.Range("J4").Value = "1003.8 hPa"
$ws.Range("K4").Value = "10.4 MJ/m2"
$ws.Range("E5").Value = "2026-02-07 18:48:01"
$ws.Range("H5").Value = "63%"
$ws.Range("J5").Value = "1003.7 hPa"
$ws.Range("K5").Value = "12.1 MJ/m2"
$ws.Range("O5").Value = "10.8 °C"
$ws.Range("E6").Value = "2026-02-07 18:48:03"
$ws.Range("J6").Value = "1005.3 hPa"
$ws.Range("K6").Value = "12.3 MJ/m2"
$ws.Range("O6").Value = "13.2 °C"
$ws.Range("E7").Value = "2026-02-07 18:48:06"
$ws.Range("J7").Value = "1004.9 hPa"
$ws.Range("K7").Value = "12.8 MJ/m2"
$ws.Range("O7").Value = "9.4 °C"
$ws.Range("E8").Value = "2026-02-07 18:48:08"
$ws.Range("K8").Value = "11.7 MJ/m2"
$ws.Range("O8").Value = "8.9 °C"
$ws.Range("E9").Value = "2026-02-07 18:48:11"
$ws.Range("H9").Value = "85%"
$ws.Range("O9").Value = "3.8 °C"
$ws.Range("E10").Value = "2026-02-07 18:48:13"
$ws.Range("O10").Value = "10.6 °C"
$ws.Range("E11").Value = "2026-02-07 18:48:16"
$ws.Range("H11").Value = "82%"
$ws.Range("J11").Value = "1006.2 hPa"
$ws.Range("K11").Value = "8.8 MJ/m2"
$ws.Range("E12").Value = "2026-02-07 18:48:19"
$ws.Range("H12").Value = "56%"
$ws.Range("K12").Value = "12.4 MJ/m2"
$ws.Range("E13").Value = "2026-02-07 18:48:22"
$ws.Range("O13").Value = "11.3 °C"
$ws.Range("E14").Value = "2026-02-07 18:48:24"
$ws.Range("H14").Value = "57%"
$ws.Range("K14").Value = "8.6 MJ/m2"
$ws.Range("L14").Value = "49.7 km/h - 249º 16:23 TU"
$ws.Range("O14").Value = "-5.5 °C"
$ws.Range("E15").Value = "2026-02-07 18:48:27"
$ws.Range("H15").Value = "69%"
$ws.Range("J15").Value = "1004.0 hPa"
$ws.Range("K15").Value = "11.6 MJ/m2"
$ws.Range("O15").Value = "10.0 °C"
$ws.Range("E16").Value = "2026-02-07 18:48:30"
$ws.Range("H16").Value = "85%"
$ws.Range("K16").Value = "6.8 MJ/m2"
$ws.Range("O16").Value = "4.1 °C"
$ws.Range("E17").Value = "2026-02-07 18:48:33"
$ws.Range("H17").Value = "83%"
$ws.Range("J17").Value = "1005.9 hPa"
$ws.Range("K17").Value = "8.8 MJ/m2"
$ws.Range("O17").Value = "5.1 °C"
$ws.Range("E18").Value = "2026-02-07 18:48:35"
$ws.Range("K18").Value = "6.0 MJ/m2"
$ws.Range("L18").Value = "19.1 km/h - 229º 17:50 TU"
$ws.Range("O18").Value = "-5.9 °C"
$ws.Range("E19").Value = "2026-02-07 18:48:37"
$ws.Range("H19").Value = "79%"
$ws.Range("J19").Value = "1007.4 hPa"
$ws.Range("K19").Value = "11.9 MJ/m2"
$ws.Range("O19").Value = "7.4 °C"
$ws.Range("E20").Value = "2026-02-07 18:48:40"
$ws.Range("K20").Value = "10.5 MJ/m2"
$ws.Range("E21").Value = "2026-02-07 18:48:43"
$ws.Range("H21").Value = "65%"
$ws.Range("J21").Value = "1003.9 hPa"
$ws.Range("K21").Value = "11.6 MJ/m2"
$ws.Range("O21").Value = "8.8 °C"
$ws.Range("E22").Value = "2026-02-07 18:48:45"
$ws.Range("K22").Value = "12.9 MJ/m2"
$ws.Range("O22").Value = "10.5 °C"
$ws.Range("E23").Value = "2026-02-07 18:48:47"
$ws.Range("H23").Value = "77%"
$ws.Range("J23").Value = "1003.7 hPa"
$ws.Range("O23").Value = "10.4 °C"
$ws.Range("E24").Value = "2026-02-07 18:48:50"
$ws.Range("H24").Value = "70%"
$ws.Range("J24").Value = "1003.2 hPa"
$ws.Range("K24").Value = "10.0 MJ/m2"
$ws.Range("E25").Value = "2026-02-07 18:48:53"
$ws.Range("H25").Value = "84%"
$ws.Range("J25").Value = "1006.0 hPa"
$ws.Range("K25").Value = "7.2 MJ/m2"
$ws.Range("L25").Value = "16.2 km/h - 204º 16:00 TU"
$ws.Range("O25").Value = "2.4 °C"
$ws.Range("E26").Value = "2026-02-07 18:48:56"
$ws.Range("K26").Value = "10.6 MJ/m2"
$ws.Range("N26").Value = "-6.2 °C 18:24 TU"
$ws.Range("O26").Value = "-2.1 °C"
$ws.Range("E27").Value = "2026-02-07 18:48:58"
$ws.Range("H27").Value = "74%"
$ws.Range("J27").Value = "1003.8 hPa"
$ws.Range("K27").Value = "11.2 MJ/m2"
$ws.Range("O27").Value = "11.6 °C"
$ws.Range("E28").Value = "2026-02-07 18:49:00"
$ws.Range("H28").Value = "79%"
$ws.Range("J28").Value = "1005.6 hPa"
$ws.Range("L28").Value = "27.0 km/h - 248º 16:26 TU"
$ws.Range("O28").Value = "5.0 °C"
$ws.Range("E29").Value = "2026-02-07 18:49:03"
$ws.Range("K29").Value = "12.5 MJ/m2"
$ws.Range("E30").Value = "2026-02-07 18:49:06"
$ws.Range("H30").Value = "70%"
$ws.Range("K30").Value = "13.7 MJ/m2"
$ws.Range("L30").Value = "54.7 km/h - 302º 18:28 TU"
$ws.Range("M30").Value = "-0.9 °C 18:14 TU"
$ws.Range("O30").Value = "-4.5 °C"
$ws.Range("E31").Value = "2026-02-07 18:49:08"
$ws.Range("H31").Value = "85%"
$ws.Range("J31").Value = "1006.7 hPa"
$ws.Range("O31").Value = "5.5 °C"
$ws.Range("E32").Value = "2026-02-07 18:49:11"
$ws.Range("H32").Value = "52%"
$ws.Range("J32").Value = "1006.5 hPa"
$ws.Range("K32").Value = "11.9 MJ/m2"
$ws.Range("O32").Value = "12.9 °C"
$ws.Range("E33").Value = "2026-02-07 18:49:14"
$ws.Range("O33").Value = "10.5 °C"
$ws.Range("E34").Value = "2026-02-07 18:49:16"
$ws.Range("H34").Value = "70%"
$ws.Range("K34").Value = "10.2 MJ/m2"
$ws.Range("O34").Value = "7.9 °C"
$ws.Range("E35").Value = "2026-02-07 18:49:19"
$ws.Range("K35").Value = "6.0 MJ/m2"
$ws.Range("E36").Value = "2026-02-07 18:49:21"
$ws.Range("H36").Value = "75%"
$ws.Range("J36").Value = "1007.8 hPa"
$ws.Range("K36").Value = "11.6 MJ/m2"
$ws.Range("O36").Value = "9.0 °C"
